# Daily attendance processing - 2026-01-15 10:39:10
# Normalizes the "Recorded By" column (G) so that entries listing both the
# automated "System" actor and a human recorder show the human first,
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value.ToString().StartsWith("System, ")) {
        $rest = $value.ToString().Substring(8)
        $cell.Value2 = $rest + ", System"
    }
}
